$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 5: Compilation success -> "no", with a note explaining why
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Missing import"

# Row 6: Runtime without error -> clear the yes/no value (keep style)
$ws.Range("B6").Value = ""

# Row 7: Assertion validity -> clear the yes/no value and its note (keep styles)
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""

# Move the active selection to B6 (matches the updated sheetView selection)
$ws.Range("B6").Select()
